# Generate Report for Handback
# Update the "generated at" timestamp text values on the Overview, zh-cn and
# de-de sheets to reflect a newly-run report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the
# 0b40fa6b-9c09-44e0-9a94-d5c1f102367b.md row (row 3, column G).
$wsOverview.Range("G3").Value = "2016-08-13 00:56:46"

# zh-cn sheet, row 3 (0b40fa6b-9c09-44e0-9a94-d5c1f102367b handback row):
#   H3 = Correspond Handoff Datetime
#   K3 = Correspond Handback DateTime
$wsZhCn.Range("H3").Value = "2016-08-13 00:56:39"
$wsZhCn.Range("K3").Value = "2016-08-13 00:57:13"

# de-de sheet, row 3 (0b40fa6b-9c09-44e0-9a94-d5c1f102367b handback row):
#   H3 = Correspond Handoff Datetime (same text as Overview!G3)
#   K3 = Correspond Handback DateTime
$wsDeDe.Range("H3").Value = "2016-08-13 00:56:46"
$wsDeDe.Range("K3").Value = "2016-08-13 00:57:23"
